# Update the "Pais" worksheet: reorder several country rows and refresh
# their statistics to the 23:16 snapshot (commit: "Update countries & provincias Spain").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp footer cell
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 23:16"

# Row edits: each entry is Row, then a hashtable of Column -> new value
$rowEdits = @(
  @{ Row=4; Cells=@{ "C"=0; "G"=0 } },
  @{ Row=9; Cells=@{ "B"=18876; "C"=5087; "D"=147; "E"=18492; "G"=30; "H"=237 } },
  @{ Row=10; Cells=@{ "D"=1587; "E"=10575; "F"=1297 } },
  @{ Row=79; Cells=@{ "A"="Marruecos"; "B"=79; "C"=16; "D"=2; "E"=74; "F"=1; "G"=1; "H"=3 } },
  @{ Row=80; Cells=@{ "A"="Brunei"; "B"=78; "C"=5; "D"=1; "E"=77; "F"=2; "G"=0; "H"=0 } },
  @{ Row=115; Cells=@{ "A"="Bolivia"; "C"=1 } },
  @{ Row=116; Cells=@{ "A"="Ghana"; "C"=5 } },
  @{ Row=121; Cells=@{ "A"="Puerto Rico"; "C"=8 } },
  @{ Row=122; Cells=@{ "A"="Montenegro"; "C"=1 } },
  @{ Row=123; Cells=@{ "A"="Guam"; "C"=2 } },
  @{ Row=132; Cells=@{ "A"="Trinidad yTobago"; "C"=0 } },
  @{ Row=133; Cells=@{ "A"="Togo"; "C"=8 } },
  @{ Row=136; Cells=@{ "A"="Kenia"; "C"=0 } },
  @{ Row=137; Cells=@{ "A"="Seychelles"; "C"=1 } },
  @{ Row=138; Cells=@{ "A"="Mongolia" } },
  @{ Row=139; Cells=@{ "A"="Mayotte"; "C"=2 } },
  @{ Row=140; Cells=@{ "A"="Tanzania"; "C"=0 } },
  @{ Row=141; Cells=@{ "A"="Kirguistan"; "C"=3 } },
  @{ Row=145; Cells=@{ "A"="Surinam"; "C"=3 } },
  @{ Row=147; Cells=@{ "A"="Bahamas"; "C"=1 } },
  @{ Row=148; Cells=@{ "A"="Gabon"; "B"=4; "C"=1; "G"=1; "H"=1 } },
  @{ Row=149; Cells=@{ "A"="Congo" } },
  @{ Row=150; Cells=@{ "A"="Madagascar"; "C"=3 } },
  @{ Row=151; Cells=@{ "A"="San Bartolome" } },
  @{ Row=152; Cells=@{ "A"="Republica de Africa Central"; "C"=2 } },
  @{ Row=153; Cells=@{ "A"="Namibia"; "C"=0 } },
  @{ Row=154; Cells=@{ "A"="Islas Virgenes de los Estados Unidos"; "E"=3; "H"=0 } },
  @{ Row=155; Cells=@{ "A"="Curazao"; "G"=0 } },
  @{ Row=156; Cells=@{ "A"="Islas Caimanes" } },
  @{ Row=157; Cells=@{ "A"="Santa Lucia"; "C"=0 } },
  @{ Row=158; Cells=@{ "A"="Zambia" } },
  @{ Row=159; Cells=@{ "A"="Benin" } },
  @{ Row=160; Cells=@{ "A"="Nueva Caledonia"; "C"=0 } },
  @{ Row=161; Cells=@{ "A"="Nicaragua" } },
  @{ Row=162; Cells=@{ "A"="Butan"; "C"=1 } },
  @{ Row=163; Cells=@{ "A"="Mauritania"; "C"=0 } },
  @{ Row=164; Cells=@{ "A"="Haiti"; "C"=2 } },
  @{ Row=165; Cells=@{ "A"="Groenlandia" } },
  @{ Row=166; Cells=@{ "A"="Guinea"; "C"=1 } },
  @{ Row=167; Cells=@{ "A"="Liberia" } },
  @{ Row=168; Cells=@{ "A"="Isla de Man"; "C"=1 } },
  @{ Row=169; Cells=@{ "A"="Bermudas" } },
  @{ Row=171; Cells=@{ "A"="Papua Nueva Guinea"; "C"=1 } },
  @{ Row=172; Cells=@{ "A"="Republica del Chad"; "C"=0 } },
  @{ Row=173; Cells=@{ "A"="Niger"; "C"=0 } },
  @{ Row=174; Cells=@{ "A"="Suazilandia" } },
  @{ Row=175; Cells=@{ "A"="Zimbabue"; "C"=1 } },
  @{ Row=176; Cells=@{ "A"="Republica de Yibuti" } },
  @{ Row=177; Cells=@{ "A"="Fiyi" } },
  @{ Row=178; Cells=@{ "A"="Cabo Verde"; "C"=1 } },
  @{ Row=179; Cells=@{ "A"="Angola"; "C"=1 } },
  @{ Row=180; Cells=@{ "A"="San Vicente y las Granadinas" } },
  @{ Row=181; Cells=@{ "A"="Santa Sede" } },
  @{ Row=182; Cells=@{ "A"="Somalia"; "C"=0 } },
  @{ Row=183; Cells=@{ "A"="El Salvador" } },
  @{ Row=184; Cells=@{ "A"="Gambia"; "C"=0 } },
  @{ Row=185; Cells=@{ "A"="Montserrat" } },
  @{ Row=187; Cells=@{ "A"="Antigua y Barbuda" } }
)

foreach ($edit in $rowEdits) {
    $r = $edit.Row
    foreach ($col in $edit.Cells.Keys) {
        $ws.Range("$col$r").Value = $edit.Cells[$col]
    }
}

